# Daily attendance processing - 2025-11-27 05:28:01
#
# The "Recorded By" column (G) holds a comma-separated list of the
# users/processes that recorded each attendance session. For this batch,
# the first name/email in the list is rotated to the end of the list for
# a specific set of rows (the rows the daily processing job touched).
#
# Example: "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column G whose "Recorded By" list needs its first entry
# rotated to the end.
$rowsToRotate = @(
    2,3,6,10,12,13,14,15,18,19,20,21,
    22,24,26,28,29,32,36,38,39,40,41,44,
    45,46,47,48,50,52,54,55,58,62,64,65,
    66,67,70,71,72,73,74,76,78,83,84,85,
    86,87,90,92,99,101,109,110,111,112,113,116,
    118,125,127,135,136,137,138,139,142,144,151,153
)

foreach ($row in $rowsToRotate) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = [string]$cell.Value2

    $parts = $current -split ',\s*'
    if ($parts.Count -gt 1) {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ', '
        $cell.Value = $rotated
    }
}
